$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / label updates -------------------------------------------------

# Timestamp footer cell
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 14:22"

# The underlying provincia data was refreshed and re-sorted by "Casos
# totales" (col B) descending; a few province names shifted rows as part of
# that re-sort.
$ws.Range("A21").Value = "Sevilla"
$ws.Range("A22").Value = "Gipuzkoa/Guipuzcoa"

$ws.Range("A26").Value = "Granada"
$ws.Range("A27").Value = "Cantabria"
$ws.Range("A28").Value = "Caceres"

# --- Numeric updates (Casos totales / Casos activos / Recuperados / Muertes)

$ws.Range("B20").Value = 2217
$ws.Range("C20").Value = 583
$ws.Range("D20").Value = 1451

$ws.Range("B21").Value = 2066
$ws.Range("C21").Value = 275
$ws.Range("D21").Value = 1614
$ws.Range("E21").Value = 177

$ws.Range("B22").Value = 2051
$ws.Range("C22").Value = 5026
$ws.Range("D22").Value = 5161
$ws.Range("E22").Value = 134

$ws.Range("B26").Value = 1832
$ws.Range("C26").Value = 360
$ws.Range("D26").Value = 1302
$ws.Range("E26").Value = 170

$ws.Range("B27").Value = 1777
$ws.Range("C27").Value = 317
$ws.Range("D27").Value = 1343
$ws.Range("E27").Value = 117

$ws.Range("B28").Value = 1776
$ws.Range("C28").Value = 237
$ws.Range("D28").Value = 1276
$ws.Range("E28").Value = 263

$ws.Range("B34").Value = 1169
$ws.Range("C34").Value = 197
$ws.Range("D34").Value = 916

$ws.Range("B35").Value = 1169
$ws.Range("C35").Value = 205
$ws.Range("D35").Value = 857
$ws.Range("E35").Value = 107

$ws.Range("B39").Value = 990
$ws.Range("C39").Value = 200
$ws.Range("D39").Value = 730
$ws.Range("E39").Value = 60

$ws.Range("B51").Value = 413
$ws.Range("C51").Value = 86
$ws.Range("D51").Value = 292
$ws.Range("E51").Value = 35

$ws.Range("B52").Value = 331
$ws.Range("D52").Value = 243
$ws.Range("E52").Value = 27
